$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "TASK0-Your responses.downloadlo"

# Update the "Marking Status" column values from "Marked" to "Posted"
$ws.Range("K2:K4").Value = "Posted"

# Convert G2:G4 from formatted/text values into plain numeric values
$ws.Range("G2").Value = 1.357
$ws.Range("G3").Value = 0.861
$ws.Range("G4").Value = 0.587

# Remove the custom numeric format previously applied to G2 (back to General)
$ws.Range("G2").ClearFormats()
